$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.874.14"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.72%  "

$ws.Range("D3").Value = "'2.460.71"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'488.60"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.80%  "

$ws.Range("D6").Value = "'152.39"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.20%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  +1.94%  "

$ws.Range("D9").Value = "'2.469.96"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.38%  "

$ws.Range("D10").Value = "'0.0996"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.58%  "

$ws.Range("D11").Value = "'5.69"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.45%  "

$ws.Range("D12").Value = "'0.334"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.19%  "

$ws.Range("E13").Value = "  +1.27%  "

$ws.Range("D14").Value = "'2.900.86"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.47%  "

$ws.Range("D15").Value = "'57.147.36"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.31%  "

$ws.Range("E16").Value = "  +2.48%  "

$ws.Range("E17").Value = "  +2.01%  "

$ws.Range("D18").Value = "'2.465.20"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("E19").Value = "  +4.65%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'322.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.15%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'10.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.02%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("E23").Value = "  +2.85%  "

$ws.Range("D24").Value = "'58.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("D25").Value = "'0.406"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.79%  "

$ws.Range("D27").Value = "'0.163"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.44%  "

$ws.Range("D28").Value = "'2.577.45"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.40%  "

$ws.Range("D29").Value = "'7.56"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.99%  "

$ws.Range("E30").Value = "  +4.46%  "

$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("D32").Value = "'150.57"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.93%  "

$ws.Range("D33").Value = "'18.29"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("D34").Value = "'1.51"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.16%  "

$ws.Range("D35").Value = "'5.21"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.82%  "

$ws.Range("E36").Value = "  +2.48%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'3.77"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.13%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'0.887"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.56%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.39"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.21%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'34.22"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.62%  "

$ws.Range("E41").Value = "  +2.63%  "

$ws.Range("D42").Value = "'0.997"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("E43").Value = "  +2.35%  "

$ws.Range("E44").Value = "  +0.89%  "

$ws.Range("D45").Value = "'0.0953"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.51%  "

$ws.Range("D46").Value = "'4.86"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.92%  "

$ws.Range("D47").Value = "'261.71"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.76%  "

$ws.Range("D48").Value = "'10.24"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.80%  "

$ws.Range("E49").Value = "  +2.69%  "

$ws.Range("D50").Value = "'17.76"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.32%  "

$ws.Range("E51").Value = "  +25.26%  "
